# Issue #21 don't add directory to playlist
# - Mark issues 20/21/22 (rows 21-23) as Status "DONE" in column C
# - Append a root-cause note to the "Playlist save not working" issue text (row 23, col H)
# - Row 23 grows to a second wrapped line, so bump its height like the other
#   wrapped rows on this sheet
# - Restore the sheet's active selection to E22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

$ws.Range("C21").Value = "DONE"
$ws.Range("C22").Value = "DONE"
$ws.Range("C23").Value = "DONE"

$ws.Range("H23").Value = "Playlist save not working.  This was caused by Bitdefender Safe Files"

$ws.Rows.Item(23).RowHeight = 29

$ws.Range("E22").Select()
